$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.031.36"
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").Value = "'1.820.10"
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  -1.28%  '
$ws.Range("D5").Value = "'310.93"
$ws.Range("E5").Value = '  -2.73%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("D7").Value = "'0.4220"
$ws.Range("E7").Value = '  -2.19%  '
$ws.Range("D8").Value = "'0.3678"
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("D9").Value = "'0.07208"
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").Value = "'0.8397"
$ws.Range("E10").Value = '  -4.42%  '
$ws.Range("D11").Value = "'20.80"
$ws.Range("E11").Value = '  -3.98%  '
$ws.Range("D12").Value = "'1.820.52"
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").Value = "'6.653"
$ws.Range("D14").Value = "'0.07052"
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = "'5.276"
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("D16").Value = "'90.17"
$ws.Range("E16").Value = '  +1.27%  '
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = '  -1.34%  '
$ws.Range("D18").Value = "'0.000008776"
$ws.Range("E18").Value = '  -2.55%  '
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("D20").Value = "'14.89"
$ws.Range("E20").Value = '  -4.12%  '
$ws.Range("D21").Value = "'27.119.15"
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").Value = "'5.125"
$ws.Range("E22").Value = '  -1.90%  '
$ws.Range("D23").Value = "'10.84"
$ws.Range("E23").Value = '  -2.39%  '
$ws.Range("D24").Value = "'2.044.60"
$ws.Range("E24").Value = '  -1.79%  '
$ws.Range("D25").Value = "'1.975"
$ws.Range("E25").Value = '  -1.01%  '
$ws.Range("D26").Value = "'151.89"
$ws.Range("E26").Value = '  -2.32%  '
$ws.Range("D27").Value = "'2.232"
$ws.Range("E27").Value = '  +2.18%  '
$ws.Range("D28").Value = "'18.25"
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("D29").Value = "'5.263"
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("D30").Value = "'116.10"
$ws.Range("E30").Value = '  -2.56%  '
$ws.Range("D31").Value = "'0.08747"
$ws.Range("E31").Value = '  -2.22%  '
$ws.Range("D32").Value = "'1.175"
$ws.Range("E32").Value = '  -4.73%  '
$ws.Range("D33").Value = "'0.7387"
$ws.Range("E33").Value = '  -5.25%  '
$ws.Range("D34").Value = "'2.914"
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").Value = "'4.416"
$ws.Range("E35").Value = '  -3.37%  '
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = '  -1.25%  '
$ws.Range("E37").Value = '  -4.03%  '
$ws.Range("D38").Value = "'0.01951"
$ws.Range("E38").Value = '  -1.39%  '
$ws.Range("D39").Value = "'0.05246"
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("D40").Value = "'7.342"
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").Value = "'2.880"
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("D42").Value = "'0.1688"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = "'0.5033"
$ws.Range("D44").Value = "'8.572"
$ws.Range("E44").Value = '  -3.00%  '
$ws.Range("D45").Value = "'10.45"
$ws.Range("E45").Value = '  -3.05%  '
$ws.Range("D46").Value = "'106.21"
$ws.Range("E46").Value = '  -2.29%  '
$ws.Range("D47").Value = "'0.4708"
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("D48").Value = "'0.9999"
$ws.Range("E48").Value = '  -1.29%  '
$ws.Range("D49").Value = "'0.06347"
$ws.Range("E49").Value = '  -2.07%  '
$ws.Range("D50").Value = "'1.880"
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").Value = "'1.646"
$ws.Range("E51").Value = '  -2.96%  '
